$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New user rows to append (id, uin, name, email, mobile)
$newRows = @(
  @(110021, 7316931025, "Magdalena Weber", "magdalena.weber@xyz.com", 932122450),
  @(110022, 9137847236, "Adrienne Hoffman", "adrienne.hoffman@xyz.com", 848488000),
  @(110023, 8428758532, "Adrienne Mcgee", "adrienne.mcgee@xyz.com", 894773246),
  @(110024, 9804209494, "Amare Coleman", "amare.coleman@xyz.com", 956554588),
  @(110025, 7105248214, "Dawson Ibarra", "dawson.ibarra@xyz.com", 765455583),
  @(110026, 9316557128, "Elvis Mcmillan", "elvis.mcmillan@xyz.com", 884282274),
  @(110027, 8103486949, "Steve George", "steve.george@xyz.com", 971073663),
  @(110028, 9601932866, "Colton Elliott", "colton.elliott@xyz.com", 809908673),
  @(110029, 9317596765, "Carolyn Rodriguez", "carolyn.rodriguez@xyz.com", 818876429)
)

$startRow = 22

# First, copy formatting (and placeholder values) for every new row down from
# the last template row, so number formats / bool style (col I) etc. match.
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $targetRow = $startRow + $i
    $srcRange = $ws.Range("A21:K21")
    $dstRange = $ws.Range("A" + $targetRow + ":K" + $targetRow)
    $srcRange.Copy($dstRange)
}

# Populate columns in field-major order (all names, then all emails, ...)
# so new shared-string entries are appended in the same order the
# original generator used (name block, then email block).
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $targetRow = $startRow + $i
    $ws.Range("A" + $targetRow).Value2 = $newRows[$i][0]
}
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $targetRow = $startRow + $i
    $ws.Range("B" + $targetRow).Value2 = $newRows[$i][1]
}
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $targetRow = $startRow + $i
    $ws.Range("C" + $targetRow).Value2 = $newRows[$i][2]
}
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $targetRow = $startRow + $i
    $ws.Range("D" + $targetRow).Value2 = $newRows[$i][3]
}
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $targetRow = $startRow + $i
    $ws.Range("E" + $targetRow).Value2 = $newRows[$i][4]
}

# Update view: scroll/selection to match the newly added block
$ws.Range("A22:K30").Select()

Write-Host "Added $($newRows.Count) rows"
